# Add carjacking data for 2022-07-25 (extends "through July 16" -> "through July 17")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the "through" column header text
$ws.Name = "Through 2022-07-17"
$ws.Range("B1").Value = "July 2022 (through July 17)"

# Updated / new cell values, keyed by neighborhood row
# Englewood (row 3)
$ws.Range("B3").Value = 8

# Auburn Gresham (row 4)
$ws.Range("B4").Value = 3

# Garfield Park (row 5)
$ws.Range("B5").Value = 5
$ws.Range("P5").Value = 8
$ws.Range("AK5").Value = 2

# Hyde Park (row 12)
$ws.Range("AK12").Value = 1

# Douglas (row 15)
$ws.Range("I15").Value = 2

# West Loop (row 18)
$ws.Range("B18").Value = 2
$ws.Range("P18").Value = 1
$ws.Range("AR18").Value = 1

# Kenwood (row 20)
$ws.Range("B20").Value = 1

# South Chicago (row 23)
$ws.Range("B23").Value = 1

# Gage Park (row 30)
$ws.Range("AR30").Value = 1

# Riverdale (row 34)
$ws.Range("I34").Value = 2

# West Town (row 38)
$ws.Range("I38").Value = 4

# New City (row 44)
$ws.Range("I44").Value = 2

# Near South Side (row 45)
$ws.Range("P45").Value = 1

# Ashburn (row 46)
$ws.Range("AD46").Value = 1

# Irving Park (row 48)
$ws.Range("AY48").Value = 1

# Calumet Heights (row 53)
$ws.Range("P53").Value = 8

# Belmont Cragin (row 56)
$ws.Range("AK56").Value = 1
$ws.Range("AY56").Value = 1
